# Fruta / hortaliza, semanal
# Insert two new weekly price rows into the Cilantro - Macroferia Regional de Talca sheet.
# New row goes in at sheet row 48 (becomes the new row 48), and another new row goes in
# at what becomes sheet row 51 after the first insert shifts everything down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole row before the current row 48 - this shifts old rows 48..70 down to 49..71.
$ws.Rows(48).Insert()

# Insert a whole row before the (new) row 51 - this shifts rows 51..71 down to 52..72.
$ws.Rows(51).Insert()

# --- Fill in the first new row (sheet row 48) ---
$ws.Cells.Item(48, 1).Value = 5
$ws.Cells.Item(48, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(48, 3).Value = "Maule"
$ws.Cells.Item(48, 4).Value = 44846
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value = 7
$ws.Cells.Item(48, 6).Value = 100112040
$ws.Cells.Item(48, 7).Value = "Cilantro"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 200
$ws.Cells.Item(48, 11).Value = 8000
$ws.Cells.Item(48, 12).Value = 8000
$ws.Cells.Item(48, 13).Value = 8000
$ws.Cells.Item(48, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(48, 15).Value = "Región del Maule"
$ws.Cells.Item(48, 16).Value = 222
$ws.Cells.Item(48, 17).Value = 36
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# --- Fill in the second new row (sheet row 51) ---
$ws.Cells.Item(51, 1).Value = 5
$ws.Cells.Item(51, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(51, 3).Value = "Maule"
$ws.Cells.Item(51, 4).Value = 44845
$ws.Cells.Item(51, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(51, 5).Value = 7
$ws.Cells.Item(51, 6).Value = 100112040
$ws.Cells.Item(51, 7).Value = "Cilantro"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 150
$ws.Cells.Item(51, 11).Value = 7000
$ws.Cells.Item(51, 12).Value = 7000
$ws.Cells.Item(51, 13).Value = 7000
$ws.Cells.Item(51, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(51, 15).Value = "Región del Maule"
$ws.Cells.Item(51, 16).Value = 194
$ws.Cells.Item(51, 17).Value = 36
$ws.Cells.Item(51, 18).Value = "Hortaliza"
